$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "INX"
$ws.Range("C4").Value = "S&P 500 Index"
$ws.Range("D4").Value = "Trading in Progress"

# E4 looks like a number ("4141.79") but must stay a text cell, matching
# the source data (same treatment as the existing E2/E3 text cells).
# Force text via a leading apostrophe, then reset the style back to
# Normal so no stray number-format/quote-prefix style sticks to the cell.
$ws.Range("E4").Value = "'4141.79"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "+24.42  +0.59%"

$ws.Range("G4").Value = 4163.5
$ws.Range("H4").Value = 4139.39
$ws.Range("I4").Value = 4607.07
$ws.Range("J4").Value = 0.84
$ws.Range("K4").Value = 4132.94
$ws.Range("L4").Value = 4117.37
$ws.Range("M4").Value = 3698.15
$ws.Range("N4").Value = 0.007
$ws.Range("O4").Value = 1142000000
